$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2135.2192
$ws.Range("J17").Value = 2365.1719
$ws.Range("L17").Value = 7095.5157
$ws.Range("N17").Value = -7431.5157

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4797.8
$ws.Range("I86").Value = 4333.3335
$ws.Range("K86").Value = 4333.3335
$ws.Range("M86").Value = -3210.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4797.8
$ws.Range("I89").Value = 4333.3335
$ws.Range("K89").Value = 21666.6675
$ws.Range("M89").Value = -16050.6675

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1567.4615
$ws.Range("I135").Value = 1303.6364
$ws.Range("K135").Value = 11732.7276
$ws.Range("M135").Value = -9197.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 40002316
$ws.Range("I45").Value = 52632840
$ws.Range("J45").Value = 5652.6665
$ws.Range("K45").Value = 52632840
$ws.Range("L45").Value = 5652.6665
$ws.Range("M45").Value = -52632463
$ws.Range("N45").Value = -6406.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5607.125
$ws.Range("I61").Value = 4860.1577
$ws.Range("K61").Value = 4860.1577
$ws.Range("M61").Value = -4648.1577

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2001.6316
$ws.Range("I102").Value = 1946.1666
$ws.Range("K102").Value = 1946.1666
$ws.Range("M102").Value = -324.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1874.8667
$ws.Range("I122").Value = 968.2222
$ws.Range("K122").Value = 2904.6666
$ws.Range("M122").Value = -454.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2440.3958
$ws.Range("I132").Value = 1891.3556
$ws.Range("K132").Value = 5674.066800000001
$ws.Range("M132").Value = -3144.066800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5607.125
$ws.Range("I136").Value = 4860.1577
$ws.Range("K136").Value = 14580.4731
$ws.Range("M136").Value = -12030.4731

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1437.3
$ws.Range("I94").Value = 1230.5
$ws.Range("J94").Value = 1747.5
$ws.Range("K94").Value = 1230.5
$ws.Range("L94").Value = 1747.5
$ws.Range("M94").Value = -779.5
$ws.Range("N94").Value = -2649.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 18089.646
$ws.Range("I105").Value = 21543.7
$ws.Range("J105").Value = 13155.286
$ws.Range("K105").Value = 21543.7
$ws.Range("L105").Value = 13155.286
$ws.Range("M105").Value = -19796.7
$ws.Range("N105").Value = -16649.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34814.79
$ws.Range("J31").Value = 90377.664
$ws.Range("L31").Value = 90377.664
$ws.Range("N31").Value = -90967.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 34814.79
$ws.Range("J34").Value = 90377.664
$ws.Range("L34").Value = 90377.664
$ws.Range("N34").Value = -90781.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 2500
$ws.Range("J54").Value = 2500
$ws.Range("L54").Value = 2500
$ws.Range("N54").Value = -3816

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1941.4
$ws.Range("I99").Value = 1867.875
$ws.Range("J99").Value = 2025.4286
$ws.Range("K99").Value = 1867.875
$ws.Range("L99").Value = 2025.4286
$ws.Range("M99").Value = -369.875
$ws.Range("N99").Value = -5021.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1941.4
$ws.Range("I126").Value = 1867.875
$ws.Range("J126").Value = 2025.4286
$ws.Range("K126").Value = 5603.625
$ws.Range("L126").Value = 6076.2858
$ws.Range("M126").Value = -3133.625
$ws.Range("N126").Value = -11016.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6117.4707
$ws.Range("I132").Value = 6416.6665
$ws.Range("K132").Value = 19249.9995
$ws.Range("M132").Value = -16719.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3680.6667
$ws.Range("I134").Value = 2278.3
$ws.Range("J134").Value = 6485.4
$ws.Range("K134").Value = 6834.900000000001
$ws.Range("L134").Value = 19456.2
$ws.Range("M134").Value = -4299.900000000001
$ws.Range("N134").Value = -24526.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 358.33334
$ws.Range("I11").Value = 950
$ws.Range("J11").Value = 62.5
$ws.Range("K11").Value = 2850
$ws.Range("L11").Value = 187.5
$ws.Range("M11").Value = -2710
$ws.Range("N11").Value = -467.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5710.3
$ws.Range("I56").Value = 5710.3
$ws.Range("K56").Value = 5710.3
$ws.Range("M56").Value = -5180.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20835464
$ws.Range("I131").Value = 25002922
$ws.Range("J131").Value = 16668007
$ws.Range("K131").Value = 75008766
$ws.Range("L131").Value = 50004021
$ws.Range("M131").Value = -75003726
$ws.Range("N131").Value = -50014101

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5460.2666
$ws.Range("I132").Value = 4380
$ws.Range("K132").Value = 39420
$ws.Range("M132").Value = -36890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4231.778
$ws.Range("I139").Value = 1651.1428
$ws.Range("J139").Value = 13264
$ws.Range("K139").Value = 4953.428400000001
$ws.Range("L139").Value = 39792
$ws.Range("M139").Value = 186.5715999999993
$ws.Range("N139").Value = -50072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2404.762
$ws.Range("I140").Value = 1583.9445
$ws.Range("K140").Value = 4751.833500000001
$ws.Range("M140").Value = 428.1664999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 22014.5
$ws.Range("J52").Value = 32999
$ws.Range("L52").Value = 32999
$ws.Range("N52").Value = -33517

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3553.5356
$ws.Range("I122").Value = 3172.3635
$ws.Range("K122").Value = 9517.0905
$ws.Range("M122").Value = -7067.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4670.8667
$ws.Range("I126").Value = 3660.182
$ws.Range("J126").Value = 7450.25
$ws.Range("K126").Value = 10980.546
$ws.Range("L126").Value = 22350.75
$ws.Range("M126").Value = -8510.545999999998
$ws.Range("N126").Value = -27290.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9501
$ws.Range("I7").Value = 4648.8335
$ws.Range("K7").Value = 4648.8335
$ws.Range("M7").Value = -4536.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7463.893
$ws.Range("I40").Value = 6946.579
$ws.Range("K40").Value = 6946.579
$ws.Range("M40").Value = -6810.579

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8950.333000000001
$ws.Range("I122").Value = 7693.5
$ws.Range("K122").Value = 23080.5
$ws.Range("M122").Value = -20630.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 9501
$ws.Range("I126").Value = 4648.8335
$ws.Range("K126").Value = 13946.5005
$ws.Range("M126").Value = -11476.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 69997.2
$ws.Range("J137").Value = 69997.2
$ws.Range("L137").Value = 69997.2
$ws.Range("N137").Value = -80197.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 17499
$ws.Range("J5").Value = 17000
$ws.Range("L5").Value = 17000
$ws.Range("N5").Value = -17224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6775.5947
$ws.Range("I122").Value = 1964.6522
$ws.Range("K122").Value = 5893.9566
$ws.Range("M122").Value = -3443.9566

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 82695
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2640.8
$ws.Range("J126").Value = 3458.7
$ws.Range("L126").Value = 10376.1
$ws.Range("N126").Value = -15316.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5455.1025
$ws.Range("I132").Value = 4826.222
$ws.Range("K132").Value = 14478.666
$ws.Range("M132").Value = -11948.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10993.571
$ws.Range("J136").Value = 15752.5
$ws.Range("L136").Value = 47257.5
$ws.Range("N136").Value = -52357.5
